$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values for Q2 and T2
$ws.Range("Q2").Value = 1
$ws.Range("T2").Value = 2

# Update row 4 values
$ws.Range("N4").Value = 0.6783525101020478
$ws.Range("P4").Value = 0.3567050202040956
$ws.Range("Q4").Value = 0.5754920420347929
$ws.Range("S4").Value = 0.1509840840695857
$ws.Range("T4").Value = 0.7796336996336996
$ws.Range("V4").Value = 0.5592673992673991

# Update B15:B21 values
$ws.Range("B15").Value = 9.85252143345858
$ws.Range("B16").Value = 5.020494655354245
$ws.Range("B17").Value = -9.351980052667205
$ws.Range("B18").Value = 6.466083189783197
$ws.Range("B19").Value = 10.07544999943667
$ws.Range("B20").Value = 13.38880928291256
$ws.Range("B21").Value = 29.12869353257992
